$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "실험없이 이벤트 효과를 추정할 수 있을까? - Difference in Differences"
$ws.Range("E3").Value = "https://lumiamitie.github.io/data/difference-in-differences/"

$ws.Range("D9").Value = "읽어볼만한 책 추천 – 딥러닝을 위한 선형대수학"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/a-book-for-read/#utm_source=rss&utm_medium=rss&utm_campaign=a-book-for-read"
